$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Review date" column (B) holds text-formatted dates such as "2011-03-09".
# Every review date whose day-of-month is "09" moves to day "10" (e.g.
# "2011-03-09" -> "2011-03-10"), while dates with any other day (e.g.
# "2024-06-02", "2024-06-26") are left untouched.

$used = $ws.UsedRange
$firstRow = $used.Row
$lastRow = $firstRow + $used.Rows.Count - 1

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 2)
    $val = $cell.Text

    if ($val -match "^\d{4}-\d{2}-09$") {
        $newVal = $val.Substring(0, $val.Length - 2) + "10"

        # Assigning a plain string like "2011-03-10" to .Value/.Formula makes
        # Excel auto-parse it as a date, turning the cell numeric and forcing
        # a new (date) number-format style onto it. Instead, enter it as a
        # text formula ("=""2011-03-10""") and then Paste Special Values-only
        # back onto itself: that freezes the literal text into the cell
        # without ever adding/altering a NumberFormat/style, exactly as the
        # original text cell had it.
        $cell.Formula = '="' + $newVal + '"'
        $cell.Copy()
        $cell.PasteSpecial(-4163)
    }
}

$excel.CutCopyMode = 0
